# Apply the "Updated symbol list" price/volume/coin-listing refresh described by
# the commit. Column D (Price) values are numeric-looking text; they are written
# with a leading apostrophe so Excel keeps storing them as text (matching the
# workbook's original inlineStr string cells) instead of coercing them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'246.86"

# Row 3
$ws.Range("D3").Value = "'26.49"

# Row 4
$ws.Range("D4").Value = "'5.094"

# Row 5
$ws.Range("D5").Value = "'0.05614"

# Row 6
$ws.Range("D6").Value = "'6.503"

# Row 8
$ws.Range("D8").Value = "'0.8439"

# Row 9
$ws.Range("B9").Value = "BitrueCoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D9").Value = "'0.02830"
$ws.Range("E9").Value = "8BitrueCoinBTR"

# Row 10
$ws.Range("B10").Value = "BitMartToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D10").Value = "'0.09394"
$ws.Range("E10").Value = "9BitMartTokenBMX"

# Row 11
$ws.Range("B11").Value = "BitForexToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D11").Value = "'0.001510"
$ws.Range("E11").Value = "10BitForexTokenBF"

# Row 12
$ws.Range("B12").Value = "One"
$ws.Range("C12").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D12").Value = "'0.0005975"
$ws.Range("E12").Value = "11OneONE"

# Row 13
$ws.Range("B13").Value = "TigerCash"
$ws.Range("C13").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D13").Value = "'0.006118"
$ws.Range("E13").Value = "12TigerCashTCH"

# Row 14
$ws.Range("B14").Value = "LEO"
$ws.Range("C14").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D14").Value = "'3.568"
$ws.Range("E14").Value = "13LEOLEO"

# Row 15
$ws.Range("B15").Value = "GateToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D15").Value = "'3.055"
$ws.Range("E15").Value = "14GateTokenGT"

# Row 16
$ws.Range("B16").Value = "BTSEToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D16").Value = "'2.118"
$ws.Range("E16").Value = "15BTSETokenBTSE"

# Row 17
$ws.Range("B17").Value = "BitpandaEcosystemToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D17").Value = "'0.3181"
$ws.Range("E17").Value = "16BitpandaEcosystemTokenBEST"

# Row 18
$ws.Range("B18").Value = "WazirX"
$ws.Range("C18").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D18").Value = "'0.1346"
$ws.Range("E18").Value = "17WazirXWRX"

# Row 19
$ws.Range("B19").Value = "MandalaExchangeToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D19").Value = "'0.06954"
$ws.Range("E19").Value = "18MandalaExchangeTokenMDX"

# Row 20
$ws.Range("B20").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C20").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D20").Value = "'0.03146"
$ws.Range("E20").Value = "19LiechtensteinCryptoassetsExchangeLCX"

# Row 21
$ws.Range("D21").Value = "'0.1319"

# Row 22
$ws.Range("D22").Value = "'3.753"

# Row 23
$ws.Range("D23").Value = "'0.04682"

# Row 25
$ws.Range("D25").Value = "'0.001248"

# Row 26
$ws.Range("D26").Value = "'0.004622"

# Row 27
$ws.Range("D27").Value = "'0.00009589"

# Row 40
$ws.Range("D40").Value = "'0.03667"

# Row 41
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1364"
$ws.Range("E41").Value = "40BKEXTokenBKKBestin24h"

# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002657"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.006143"
$ws.Range("E43").Value = "42KickTokenKICK"

# Row 44
$ws.Range("D44").Value = "'0.008685"

# Row 45
$ws.Range("D45").Value = "'0.00005290"

# Row 47
$ws.Range("D47").Value = "'0.1599"
